$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "37.269.71"
$r.Style = "Normal"
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "  +2.08%  "
$r.Style = "Normal"
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "2.001.87"
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "  +2.73%  "
$r.Style = "Normal"
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "  -0.07%  "
$r.Style = "Normal"
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "246.63"
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "  +1.64%  "
$r.Style = "Normal"
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.628"
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "  +2.88%  "
$r.Style = "Normal"
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "60.64"
$r.Style = "Normal"
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "  +5.53%  "
$r.Style = "Normal"
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "  +2.19%  "
$r.Style = "Normal"
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.0800"
$r.Style = "Normal"
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "  +2.19%  "
$r.Style = "Normal"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "  +1.42%  "
$r.Style = "Normal"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "14.99"
$r.Style = "Normal"
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "  +10.67%  "
$r.Style = "Normal"
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "22.75"
$r.Style = "Normal"
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "  +7.31%  "
$r.Style = "Normal"
$r = $ws.Range("B14")
$r.NumberFormat = "@"
$r.Value = "WrappedliquidstakedEther2.0"
$r.Style = "Normal"
$r = $ws.Range("C14")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$r.Style = "Normal"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "2.298.41"
$r.Style = "Normal"
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "  +2.70%  "
$r.Style = "Normal"
$r = $ws.Range("B15")
$r.NumberFormat = "@"
$r.Value = "Polygon"
$r.Style = "Normal"
$r = $ws.Range("C15")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$r.Style = "Normal"
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.845"
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "  +3.11%  "
$r.Style = "Normal"
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "5.44"
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "  +3.24%  "
$r.Style = "Normal"
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "2.009.62"
$r.Style = "Normal"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "  +2.99%  "
$r.Style = "Normal"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "37.205.69"
$r.Style = "Normal"
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "  +2.47%  "
$r.Style = "Normal"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "70.47"
$r.Style = "Normal"
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "  +2.14%  "
$r.Style = "Normal"
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "0.0₃0864"
$r.Style = "Normal"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "  +2.68%  "
$r.Style = "Normal"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "5.19"
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "  +3.92%  "
$r.Style = "Normal"
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "230.76"
$r.Style = "Normal"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "  +1.63%  "
$r.Style = "Normal"
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "  +0.12%  "
$r.Style = "Normal"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "2.51"
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "  +1.78%  "
$r.Style = "Normal"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "2.36"
$r.Style = "Normal"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "  +0.11%  "
$r.Style = "Normal"
$r = $ws.Range("B26")
$r.NumberFormat = "@"
$r.Value = "Kaspa"
$r.Style = "Normal"
$r = $ws.Range("C26")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$r.Style = "Normal"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "0.144"
$r.Style = "Normal"
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "  +4.96%  "
$r.Style = "Normal"
$r = $ws.Range("B27")
$r.NumberFormat = "@"
$r.Value = "Cosmos"
$r.Style = "Normal"
$r = $ws.Range("C27")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$r.Style = "Normal"
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "9.38"
$r.Style = "Normal"
$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "  +3.91%  "
$r.Style = "Normal"
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "163.75"
$r.Style = "Normal"
$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "  +2.27%  "
$r.Style = "Normal"
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "19.69"
$r.Style = "Normal"
$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = "  +2.42%  "
$r.Style = "Normal"
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "1.33"
$r.Style = "Normal"
$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = "  +15.71%  "
$r.Style = "Normal"
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.121"
$r.Style = "Normal"
$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = "  +1.56%  "
$r.Style = "Normal"
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "4.83"
$r.Style = "Normal"
$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = "  +4.06%  "
$r.Style = "Normal"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "0.0648"
$r.Style = "Normal"
$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = "  +6.65%  "
$r.Style = "Normal"
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "4.54"
$r.Style = "Normal"
$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = "  +6.28%  "
$r.Style = "Normal"
$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = "  +7.32%  "
$r.Style = "Normal"
$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = "  +0.16%  "
$r.Style = "Normal"
$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = "  +2.13%  "
$r.Style = "Normal"
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "3.29"
$r.Style = "Normal"
$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "  -2.67%  "
$r.Style = "Normal"
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "5.50"
$r.Style = "Normal"
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "  +4.70%  "
$r.Style = "Normal"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.0979"
$r.Style = "Normal"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "  +1.45%  "
$r.Style = "Normal"
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "  +0.86%  "
$r.Style = "Normal"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "1.18"
$r.Style = "Normal"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "  +2.78%  "
$r.Style = "Normal"
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.0214"
$r.Style = "Normal"
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "  +2.52%  "
$r.Style = "Normal"
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "16.53"
$r.Style = "Normal"
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "  +5.57%  "
$r.Style = "Normal"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "90.69"
$r.Style = "Normal"
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "  +4.44%  "
$r.Style = "Normal"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "1.370.71"
$r.Style = "Normal"
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "  +1.16%  "
$r.Style = "Normal"
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "  +3.00%  "
$r.Style = "Normal"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "7.26"
$r.Style = "Normal"
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "  +2.85%  "
$r.Style = "Normal"
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "  +1.01%  "
$r.Style = "Normal"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "2.01"
$r.Style = "Normal"
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "  +15.82%  "
$r.Style = "Normal"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "46.02"
$r.Style = "Normal"
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "  +6.83%  "
$r.Style = "Normal"
